$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Femacal de La Calera" - Naranja (Valencia, Primera)
# It must be inserted as a new row 1142, pushing the existing rows 1142-1224 down to 1143-1225.
# The new row reuses the same descriptive fields (A,B,C,E,F,G,H,I,J,K,L,N,Q,R,T) as the
# row that used to occupy position 1142, only the date and price-related figures change.

# Capture the descriptive values from the row that is about to be pushed down (old row 1142)
$A = $ws.Range("A1142").Value2
$B = $ws.Range("B1142").Value2
$C = $ws.Range("C1142").Value2
$E = $ws.Range("E1142").Value2
$F = $ws.Range("F1142").Value2
$G = $ws.Range("G1142").Value2
$H = $ws.Range("H1142").Value2
$I = $ws.Range("I1142").Value2
$J = $ws.Range("J1142").Value2
$K = $ws.Range("K1142").Value2
$L = $ws.Range("L1142").Value2
$N = $ws.Range("N1142").Value2
$Q = $ws.Range("Q1142").Value2
$R = $ws.Range("R1142").Value2
$T = $ws.Range("T1142").Value2

# Insert a new blank row at position 1142; rows 1142-1224 shift to 1143-1225
$ws.Rows(1142).Insert()

# Fill in the new row 1142 with the carried-over descriptive fields ...
$ws.Range("A1142").Value = $A
$ws.Range("B1142").Value = $B
$ws.Range("C1142").Value = $C
$ws.Range("D1142").Value = 45013
$ws.Range("E1142").Value = $E
$ws.Range("F1142").Value = $F
$ws.Range("G1142").Value = $G
$ws.Range("H1142").Value = $H
$ws.Range("I1142").Value = $I
$ws.Range("J1142").Value = $J
$ws.Range("K1142").Value = $K
$ws.Range("L1142").Value = $L
$ws.Range("M1142").Value = 135
$ws.Range("N1142").Value = $N
$ws.Range("O1142").Value = 8000
$ws.Range("P1142").Value = 7496
$ws.Range("Q1142").Value = $Q
$ws.Range("R1142").Value = $R
$ws.Range("S1142").Value = 577
$ws.Range("T1142").Value = $T

# Keep the date column formatted consistently with the rest of the sheet
$ws.Range("D1142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
